# ----------------------------------------------------------------------------
# Edit: Jogos_da_Semana_FlashScore_2024-11-08.xlsx
# - Fix header order for AW1:BC1 (Odd_CS_*_HT columns were shifted by one)
# - Replace the single match row (row 2) with new match data (BULGARIA)
# - Append a new match row (row 3, ISRAEL)
# ----------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Fix the AW1:BC1 header ordering (left-rotate by one column) ---------
$ws.Range("AW1").Value = "Odd_CS_0-1_HT"
$ws.Range("AX1").Value = "Odd_CS_0-2_HT"
$ws.Range("AY1").Value = "Odd_CS_1-2_HT"
$ws.Range("AZ1").Value = "Odd_CS_0-3_HT"
$ws.Range("BA1").Value = "Odd_CS_1-3_HT"
$ws.Range("BB1").Value = "Odd_CS_2-3_HT"
$ws.Range("BC1").Value = "Odd_CS_3-3_HT"

# --- 2) Row 2: overwrite with the new match (BULGARIA - PARVA LIGA) --------
# Columns B (Date) hold literal text like "08/11/2024"; assigning that string
# directly would be auto-parsed into a real date serial by the COM layer, so
# we quote-prefix it to force literal text, then copy a plain (unstyled) cell
# style over it so the number formatting introduced by the quote-prefix does
# not stick around.
$ws.Range("A2").Value = "z1aTWoV1"
$ws.Range("B2").Value = "'08/11/2024"
$ws.Range("B2").Style = $ws.Range("C2").Style
$ws.Range("C2").Value = "10:00"
$ws.Range("D2").Value = "BULGARIA - PARVA LIGA"
$ws.Range("E2").Value = "Arda"
$ws.Range("F2").Value = "Botev Vratsa"
$ws.Range("G2").Value = 1.57
$ws.Range("H2").Value = 3.75
$ws.Range("I2").Value = 6.25
$ws.Range("J2").Value = 2.2
$ws.Range("K2").Value = 2.1
$ws.Range("L2").Value = 6.5
$ws.Range("M2").Value = 1.07
$ws.Range("N2").Value = 9
$ws.Range("O2").Value = 1.36
$ws.Range("P2").Value = 3
$ws.Range("Q2").Value = 2.15
$ws.Range("R2").Value = 1.67
$ws.Range("S2").Value = 1.5
$ws.Range("T2").Value = 2.5
$ws.Range("U2").Value = 2.25
$ws.Range("V2").Value = 1.57
$ws.Range("W2").Value = 5.5
$ws.Range("X2").Value = 6.5
$ws.Range("Y2").Value = 9
$ws.Range("Z2").Value = 11
$ws.Range("AA2").Value = 15
$ws.Range("AB2").Value = 34
$ws.Range("AC2").Value = 8
$ws.Range("AD2").Value = 7.5
$ws.Range("AE2").Value = 21
$ws.Range("AF2").Value = 81
$ws.Range("AG2").Value = 1250
$ws.Range("AH2").Value = 13
$ws.Range("AI2").Value = 29
$ws.Range("AJ2").Value = 21
$ws.Range("AK2").Value = 67
$ws.Range("AL2").Value = 51
$ws.Range("AM2").Value = 51
$ws.Range("AN2").Value = 3.4
$ws.Range("AO2").Value = 8
$ws.Range("AP2").Value = 23
$ws.Range("AQ2").Value = 26
$ws.Range("AR2").Value = 51
$ws.Range("AS2").Value = 201
$ws.Range("AT2").Value = 2.5
$ws.Range("AU2").Value = 10
$ws.Range("AV2").Value = 81
$ws.Range("AW2").Value = 7.5
$ws.Range("AX2").Value = 34
$ws.Range("AY2").Value = 41
$ws.Range("AZ2").Value = 151
$ws.Range("BA2").Value = 201
$ws.Range("BB2").Value = 351
$ws.Range("BC2").Value = 51
$ws.Range("BD2").Value = 51

# --- 3) Row 3: brand-new match (ISRAEL - LIGAT HA'AL) -----------------------
$ws.Range("A3").Value = "prNuQMTF"
$ws.Range("B3").Value = "'08/11/2024"
$ws.Range("B3").Style = $ws.Range("C3").Style
$ws.Range("C3").Value = "10:00"
$ws.Range("D3").Value = "ISRAEL - LIGAT HA'AL"
$ws.Range("E3").Value = "Hapoel Jerusalem"
$ws.Range("F3").Value = "Maccabi Haifa"
$ws.Range("G3").Value = 7
$ws.Range("H3").Value = 4.5
$ws.Range("I3").Value = 1.33
$ws.Range("J3").Value = 7
$ws.Range("K3").Value = 2.6
$ws.Range("L3").Value = 1.83
$ws.Range("M3").Value = 1.03
$ws.Range("N3").Value = 17
$ws.Range("O3").Value = 1.17
$ws.Range("P3").Value = 5
$ws.Range("Q3").Value = 1.53
$ws.Range("R3").Value = 2.4
$ws.Range("S3").Value = 1.29
$ws.Range("T3").Value = 3.5
$ws.Range("U3").Value = 1.83
$ws.Range("V3").Value = 1.83
$ws.Range("W3").Value = 21
$ws.Range("X3").Value = 41
$ws.Range("Y3").Value = 23
$ws.Range("Z3").Value = 81
$ws.Range("AA3").Value = 51
$ws.Range("AB3").Value = 51
$ws.Range("AC3").Value = 15
$ws.Range("AD3").Value = 9.5
$ws.Range("AE3").Value = 19
$ws.Range("AF3").Value = 51
$ws.Range("AG3").Value = 251
$ws.Range("AH3").Value = 8.5
$ws.Range("AI3").Value = 7.5
$ws.Range("AJ3").Value = 8.5
$ws.Range("AK3").Value = 9.5
$ws.Range("AL3").Value = 11
$ws.Range("AM3").Value = 23
$ws.Range("AN3").Value = 9
$ws.Range("AO3").Value = 34
$ws.Range("AP3").Value = 34
$ws.Range("AQ3").Value = 126
$ws.Range("AR3").Value = 126
$ws.Range("AS3").Value = 201
$ws.Range("AT3").Value = 3.5
$ws.Range("AU3").Value = 8.5
$ws.Range("AV3").Value = 51
$ws.Range("AW3").Value = 3.5
$ws.Range("AX3").Value = 6.5
$ws.Range("AY3").Value = 17
$ws.Range("AZ3").Value = 17
$ws.Range("BA3").Value = 41
$ws.Range("BB3").Value = 101
$ws.Range("BC3").Value = 51
$ws.Range("BD3").Value = 51
